# About to start trial 2
# Adds a second model/trial row to the results sheet, adds a new
# "Other" notes column, and updates the best-accuracy note for trial 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New content is entered in this order so that the workbook's shared
# string table is rebuilt with the same ordering as the authored file ---

# New trial-2 model name (row 3)
$ws.Range("B3").Value = "Simple_MLP([40*(2*context_size+1), 256, 128, 71]) + Softmax"

# New "Other" header column
$ws.Range("G1").Value = "Other"

# Updated accuracy note for trial 1
$ws.Range("F2").Value = "66.68% - Epoch 11"

# New note for trial 1 in the "Other" column
$ws.Range("G2").Value = "65% Epoch 5, 66.48% Epoch 10"

# Remaining cells of the new trial-2 row
$ws.Range("A3").Value = 2
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = "Adam"
$ws.Range("E3").Value = 0.0001
$ws.Range("E3").NumberFormat = $ws.Range("E2").NumberFormat

# --- Column width adjustments to fit the new/updated content ---
$ws.Columns.Item(2).ColumnWidth = 49.5
$ws.Columns.Item(5).ColumnWidth = 13.333333333333332
$ws.Columns.Item(6).ColumnWidth = 15.833333333333332
$ws.Columns.Item(7).ColumnWidth = 29.666666666666664
